# Add a "Save" column (H) to the s_vals sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - same style as the other header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for H2:H16
$values = @(0, 0, 1, 0, 1, 0, 1, 1, 0, 0, 0, 1, 0, 0, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
